$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.040308356285095
$ws.Range("B1").Value = 0.8119840621948242
$ws.Range("C1").Value = 3.228909969329834
$ws.Range("D1").Value = 3.177035808563232
$ws.Range("E1").Value = 0.941809356212616
